$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$addresses = @(
    "B2", "C2", "D2", "E2", "F2", "I2", "J2", "K2", "L2", "M2",
    "B3", "C3", "D3", "E3", "F3", "I3", "J3", "K3", "L3", "M3",
    "B4", "C4", "D4", "E4", "F4", "I4", "J4", "K4", "L4", "M4",
    "B5", "C5", "D5", "E5", "F5", "I5", "J5", "K5", "L5", "M5",
    "B6", "C6", "D6", "E6", "F6", "I6", "J6", "K6", "L6", "M6",
    "B7", "C7", "D7", "E7", "F7", "I7", "J7", "K7", "L7", "M7",
    "B8", "C8", "D8", "E8", "F8", "I8", "J8", "K8", "L8", "M8",
    "B9", "C9", "D9", "E9", "F9", "I9", "J9", "K9", "L9", "M9",
    "B10", "C10", "D10", "E10", "F10", "I10", "J10", "K10", "L10", "M10",
    "B11", "C11", "D11", "E11", "F11", "I11", "J11", "K11", "L11", "M11",
    "B12", "C12", "D12", "E12", "F12", "I12", "J12", "K12", "L12", "M12",
    "B13", "C13", "D13", "E13", "F13", "I13", "J13", "K13", "L13", "M13",
    "B14", "C14", "D14", "E14", "F14", "I14", "J14", "K14", "L14", "M14",
    "B15", "C15", "D15", "E15", "F15", "I15", "J15", "K15", "L15", "M15",
    "B16", "C16", "D16", "E16", "F16", "I16", "J16", "K16", "L16", "M16",
    "B17", "C17", "D17", "E17", "F17", "I17", "J17", "K17", "L17", "M17",
    "B18", "C18", "D18", "E18", "F18", "I18", "J18", "K18", "L18", "M18",
    "B19", "C19", "D19", "E19", "F19", "I19", "J19", "K19", "L19", "M19",
    "B20", "C20", "D20", "E20", "F20", "I20", "J20", "K20", "L20", "M20",
    "B21", "C21", "D21", "E21", "F21", "I21", "J21", "K21", "L21", "M21",
    "B22", "C22", "D22", "E22", "F22", "I22", "J22", "K22", "L22", "M22",
    "B23", "C23", "D23", "E23", "F23", "I23", "J23", "K23", "L23", "M23",
    "B24", "C24", "D24", "E24", "F24", "I24", "J24", "K24", "L24", "M24",
    "B25", "C25", "D25", "E25", "F25", "I25", "J25", "K25", "L25", "M25"
)

$values = @(
    1.02, 1.051558125831511, 1.053748318504431, 1.048251097400144, 1.059505324490594, 1.049647275507857, 1.056585365705097, 1.056493430923029, 1.051011464325533, 1.062234645293643,
    1.02, 1.05294731746347, 1.054866382000574, 1.049473553825967, 1.0610855122101, 1.050193171645195, 1.057622546760413, 1.057423820220612, 1.052044866240877, 1.063627139126467,
    1.02, 1.053844158736299, 1.055588056272542, 1.050262215559527, 1.062106466013311, 1.05054404292949, 1.058291167769236, 1.058023473583255, 1.052710639686535, 1.064526154751455,
    1.02, 1.054220705838856, 1.055891025821862, 1.050593213435466, 1.062535316866567, 1.050690987816083, 1.058571662584958, 1.058275005659778, 1.05298984142194, 1.064903625239759,
    1.02, 1.054283901543089, 1.055941871140527, 1.050648757059594, 1.062607302070337, 1.050715627692171, 1.058618724275321, 1.058317206189722, 1.053036680380433, 1.064966976537437,
    1.02, 1.053849192072976, 1.055592106218936, 1.050266640543588, 1.06211219773551, 1.050546008615752, 1.058294918079007, 1.058026836769608, 1.052714373094769, 1.064531200389616,
    1.02, 1.052028040645142, 1.054126546829946, 1.04866472288984, 1.060039679217146, 1.049832254490903, 1.056936408167732, 1.056808354454157, 1.051361313085015, 1.062705668486253,
    1.02, 1.048802818774899, 1.051530099403937, 1.045823635147012, 1.056375490799912, 1.048556291757873, 1.054523083058827, 1.054642834580094, 1.048954496490032, 1.059473002969653,
    1.02, 1.046641335068541, 1.049789401259321, 1.043916832333014, 1.053923966951429, 1.047693160610111, 1.052900729974029, 1.053186444732701, 1.047334381436358, 1.057306724343322,
    1.02, 1.045702588659876, 1.049033270768513, 1.04308804476538, 1.052860232595178, 1.047316401768258, 1.052194951347475, 1.052552720887711, 1.046629069129724, 1.056365937556086,
    1.02, 1.045353464631602, 1.048752043138157, 1.042779717460563, 1.052464771539168, 1.047175998949094, 1.051932292012034, 1.052316855613463, 1.046366507086501, 1.05601606070102,
    1.02, 1.045428372617767, 1.048812384144973, 1.042845876514342, 1.052549614965605, 1.047206136619744, 1.051988656157945, 1.052367470983133, 1.04642285380463, 1.056091129947336,
    1.02, 1.045673738824015, 1.049010031934424, 1.043062568149238, 1.052827550716797, 1.047304805392709, 1.052173250124121, 1.052533233858449, 1.046607377499546, 1.056337025389985,
    1.02, 1.045824859530916, 1.049131760405309, 1.043196015327771, 1.052998750358128, 1.047365537681034, 1.05228691784696, 1.052635303088835, 1.046720991854296, 1.056488472837698,
    1.02, 1.046703576553772, 1.049839532056078, 1.043971769574214, 1.053994516035383, 1.047718100899351, 1.052947500234685, 1.053228437070134, 1.047381110123395, 1.05736910196938,
    1.02, 1.047254013476258, 1.050282851910468, 1.044457535885567, 1.054618534580412, 1.04793844350444, 1.053360980209706, 1.053599660387693, 1.0477941640935, 1.057920747676792,
    1.02, 1.047574803584365, 1.050541202111873, 1.044740573799725, 1.054982301712597, 1.048066674792103, 1.053601839255905, 1.053815890265079, 1.048034726175005, 1.058242246330084,
    1.02, 1.047684139142251, 1.050629253922891, 1.044837031556055, 1.055106301220023, 1.048110349153241, 1.053683912455193, 1.053889568651645, 1.048116689905065, 1.058351824083922,
    1.02, 1.047194984846378, 1.050235311806031, 1.044405448990246, 1.054551605337648, 1.047914832942194, 1.053316650547439, 1.053559862552458, 1.047749885155168, 1.057861588997981,
    1.02, 1.045601496590288, 1.048951839777127, 1.042998771147158, 1.052745715115702, 1.047275762572355, 1.052118905713414, 1.05248443390716, 1.04655305588083, 1.056264627172257,
    1.02, 1.04459710407348, 1.048142742522503, 1.042111563055202, 1.051608291689326, 1.046871302666644, 1.051362929342583, 1.051805534915627, 1.045797215718906, 1.055258079800747,
    1.02, 1.045129792081082, 1.048571864304036, 1.042582154657737, 1.052211453548915, 1.04708596730733, 1.051763964788072, 1.052165693539451, 1.046198220740171, 1.055791907487091,
    1.02, 1.047221658179142, 1.050256793846133, 1.044428985747513, 1.054581848437734, 1.047925502437677, 1.053336682192773, 1.053577846403299, 1.047769894030335, 1.057888321083478,
    1.02, 1.049638577237201, 1.052203032421537, 1.046560338398062, 1.057324268620001, 1.048888343203396, 1.055149330934433, 1.055204888483483, 1.049579430592633, 1.060310654973628
)

for ($i = 0; $i -lt $addresses.Length; $i++) {
    $ws.Range($addresses[$i]).Value = $values[$i]
}

Write-Output "Updated $($addresses.Length) cells"
